$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (92 cell changes) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(3, 8).Value = 40657  # H3
$ws.Cells.Item(3, 10).Value = 40657  # J3
$ws.Cells.Item(3, 12).Value = 40657  # L3
$ws.Cells.Item(3, 14).Value = -40885  # N3
$ws.Cells.Item(41, 8).Value = 493.04877  # H41
$ws.Cells.Item(41, 9).Value = 510.54544  # I41
$ws.Cells.Item(41, 10).Value = 472.78946  # J41
$ws.Cells.Item(41, 11).Value = 510.54544  # K41
$ws.Cells.Item(41, 12).Value = 472.78946  # L41
$ws.Cells.Item(41, 13).Value = -70.54543999999999  # M41
$ws.Cells.Item(41, 14).Value = -1352.78946  # N41
$ws.Cells.Item(70, 8).Value = 3950.261  # H70
$ws.Cells.Item(70, 9).Value = 3453.3125  # I70
$ws.Cells.Item(70, 10).Value = 5086.143  # J70
$ws.Cells.Item(70, 11).Value = 10359.9375  # K70
$ws.Cells.Item(70, 12).Value = 15258.429  # L70
$ws.Cells.Item(70, 13).Value = -10089.9375  # M70
$ws.Cells.Item(70, 14).Value = -15798.429  # N70
$ws.Cells.Item(73, 8).Value = 3950.261  # H73
$ws.Cells.Item(73, 9).Value = 3453.3125  # I73
$ws.Cells.Item(73, 10).Value = 5086.143  # J73
$ws.Cells.Item(73, 11).Value = 10359.9375  # K73
$ws.Cells.Item(73, 12).Value = 15258.429  # L73
$ws.Cells.Item(73, 13).Value = -9423.9375  # M73
$ws.Cells.Item(73, 14).Value = -17130.429  # N73
$ws.Cells.Item(74, 8).Value = 4880.769  # H74
$ws.Cells.Item(74, 10).Value = 4880.769  # J74
$ws.Cells.Item(74, 12).Value = 4880.769  # L74
$ws.Cells.Item(74, 14).Value = -6752.769  # N74
$ws.Cells.Item(77, 8).Value = 4880.769  # H77
$ws.Cells.Item(77, 10).Value = 4880.769  # J77
$ws.Cells.Item(77, 12).Value = 24403.845  # L77
$ws.Cells.Item(77, 14).Value = -33763.845  # N77
$ws.Cells.Item(80, 8).Value = 697.61536  # H80
$ws.Cells.Item(80, 9).Value = 761.75  # I80
$ws.Cells.Item(80, 10).Value = 669.1111  # J80
$ws.Cells.Item(80, 11).Value = 2285.25  # K80
$ws.Cells.Item(80, 12).Value = 2007.3333  # L80
$ws.Cells.Item(80, 13).Value = -1287.25  # M80
$ws.Cells.Item(80, 14).Value = -4003.3333  # N80
$ws.Cells.Item(83, 8).Value = 697.61536  # H83
$ws.Cells.Item(83, 9).Value = 761.75  # I83
$ws.Cells.Item(83, 10).Value = 669.1111  # J83
$ws.Cells.Item(83, 11).Value = 6855.75  # K83
$ws.Cells.Item(83, 12).Value = 6021.9999  # L83
$ws.Cells.Item(83, 13).Value = -1863.75  # M83
$ws.Cells.Item(83, 14).Value = -16005.9999  # N83
$ws.Cells.Item(86, 8).Value = 10002616  # H86
$ws.Cells.Item(86, 9).Value = 16669129  # I86
$ws.Cells.Item(86, 10).Value = 2847.5  # J86
$ws.Cells.Item(86, 11).Value = 16669129  # K86
$ws.Cells.Item(86, 12).Value = 2847.5  # L86
$ws.Cells.Item(86, 13).Value = -16668006  # M86
$ws.Cells.Item(86, 14).Value = -5093.5  # N86
$ws.Cells.Item(89, 8).Value = 10002616  # H89
$ws.Cells.Item(89, 9).Value = 16669129  # I89
$ws.Cells.Item(89, 10).Value = 2847.5  # J89
$ws.Cells.Item(89, 11).Value = 83345645  # K89
$ws.Cells.Item(89, 12).Value = 14237.5  # L89
$ws.Cells.Item(89, 13).Value = -83340029  # M89
$ws.Cells.Item(89, 14).Value = -25469.5  # N89
$ws.Cells.Item(92, 8).Value = 768.38464  # H92
$ws.Cells.Item(92, 9).Value = 768.38464  # I92
$ws.Cells.Item(92, 11).Value = 768.38464  # K92
$ws.Cells.Item(92, 13).Value = 479.61536  # M92
$ws.Cells.Item(102, 8).Value = 40657  # H102
$ws.Cells.Item(102, 10).Value = 40657  # J102
$ws.Cells.Item(102, 12).Value = 40657  # L102
$ws.Cells.Item(102, 14).Value = -47147  # N102
$ws.Cells.Item(112, 8).Value = 2383.125  # H112
$ws.Cells.Item(112, 10).Value = 3285  # J112
$ws.Cells.Item(112, 12).Value = 9855  # L112
$ws.Cells.Item(112, 14).Value = -12071  # N112
$ws.Cells.Item(115, 8).Value = 1564.8  # H115
$ws.Cells.Item(115, 9).Value = 1702  # I115
$ws.Cells.Item(115, 11).Value = 5106  # K115
$ws.Cells.Item(115, 13).Value = -3539  # M115
$ws.Cells.Item(116, 8).Value = 7671.5  # H116
$ws.Cells.Item(116, 9).Value = 4509.6665  # I116
$ws.Cells.Item(116, 11).Value = 4509.6665  # K116
$ws.Cells.Item(116, 13).Value = -1067.6665  # M116
$ws.Cells.Item(137, 8).Value = 57500.723  # H137
$ws.Cells.Item(137, 9).Value = 1794.2667  # I137
$ws.Cells.Item(137, 11).Value = 5382.800099999999  # K137
$ws.Cells.Item(137, 13).Value = -2832.800099999999  # M137
$ws.Cells.Item(138, 8).Value = 3466.4546  # H138
$ws.Cells.Item(138, 9).Value = 2081.875  # I138
$ws.Cells.Item(138, 10).Value = 4769.5884  # J138
$ws.Cells.Item(138, 11).Value = 6245.625  # K138
$ws.Cells.Item(138, 12).Value = 14308.7652  # L138
$ws.Cells.Item(138, 13).Value = -1105.625  # M138
$ws.Cells.Item(138, 14).Value = -24588.7652  # N138

# ---- Sheet: ARM (60 cell changes) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 147283.08  # H32
$ws.Cells.Item(32, 9).Value = 137673.73  # I32
$ws.Cells.Item(32, 11).Value = 137673.73  # K32
$ws.Cells.Item(32, 13).Value = -137386.73  # M32
$ws.Cells.Item(37, 8).Value = 16534  # H37
$ws.Cells.Item(37, 10).Value = 0  # J37
$ws.Cells.Item(37, 12).Value = 0  # L37
$ws.Cells.Item(37, 14).ClearContents()  # N37: was -30546
$ws.Cells.Item(38, 8).Value = 5000  # H38
$ws.Cells.Item(38, 9).Value = 5000  # I38
$ws.Cells.Item(38, 10).Value = 0  # J38
$ws.Cells.Item(38, 11).Value = 5000  # K38
$ws.Cells.Item(38, 12).Value = 0  # L38
$ws.Cells.Item(38, 13).Value = -4533  # M38
$ws.Cells.Item(38, 14).ClearContents()  # N38: was -20934
$ws.Cells.Item(39, 8).Value = 6495.75  # H39
$ws.Cells.Item(39, 10).Value = 12500  # J39
$ws.Cells.Item(39, 12).Value = 12500  # L39
$ws.Cells.Item(39, 14).Value = -13540  # N39
$ws.Cells.Item(61, 8).Value = 3094.6  # H61
$ws.Cells.Item(61, 9).Value = 3160.6667  # I61
$ws.Cells.Item(61, 11).Value = 3160.6667  # K61
$ws.Cells.Item(61, 13).Value = -2948.6667  # M61
$ws.Cells.Item(63, 8).Value = 0  # H63
$ws.Cells.Item(63, 9).Value = 0  # I63
$ws.Cells.Item(63, 10).Value = 0  # J63
$ws.Cells.Item(63, 11).Value = 0  # K63
$ws.Cells.Item(63, 12).Value = 0  # L63
$ws.Cells.Item(63, 13).ClearContents()  # M63: was -14
$ws.Cells.Item(63, 14).ClearContents()  # N63: was -2372
$ws.Cells.Item(66, 8).Value = 0  # H66
$ws.Cells.Item(66, 9).Value = 0  # I66
$ws.Cells.Item(66, 10).Value = 0  # J66
$ws.Cells.Item(66, 11).Value = 0  # K66
$ws.Cells.Item(66, 12).Value = 0  # L66
$ws.Cells.Item(66, 13).ClearContents()  # M66: was -68
$ws.Cells.Item(66, 14).ClearContents()  # N66: was -11864
$ws.Cells.Item(102, 8).Value = 2426.1177  # H102
$ws.Cells.Item(102, 9).Value = 2323.4666  # I102
$ws.Cells.Item(102, 11).Value = 2323.4666  # K102
$ws.Cells.Item(102, 13).Value = -701.4666000000002  # M102
$ws.Cells.Item(112, 8).Value = 70386  # H112
$ws.Cells.Item(112, 10).Value = 70386  # J112
$ws.Cells.Item(112, 12).Value = 70386  # L112
$ws.Cells.Item(112, 14).Value = -73340  # N112
$ws.Cells.Item(122, 8).Value = 30809.143  # H122
$ws.Cells.Item(122, 9).Value = 41610  # I122
$ws.Cells.Item(122, 10).Value = 3807  # J122
$ws.Cells.Item(122, 11).Value = 124830  # K122
$ws.Cells.Item(122, 12).Value = 11421  # L122
$ws.Cells.Item(122, 13).Value = -122380  # M122
$ws.Cells.Item(122, 14).Value = -16321  # N122
$ws.Cells.Item(132, 8).Value = 3027.4285  # H132
$ws.Cells.Item(132, 9).Value = 2838.4  # I132
$ws.Cells.Item(132, 11).Value = 8515.200000000001  # K132
$ws.Cells.Item(132, 13).Value = -5985.200000000001  # M132
$ws.Cells.Item(136, 8).Value = 3094.6  # H136
$ws.Cells.Item(136, 9).Value = 3160.6667  # I136
$ws.Cells.Item(136, 11).Value = 9482.000100000001  # K136
$ws.Cells.Item(136, 13).Value = -6932.000100000001  # M136

# ---- Sheet: BSM (59 cell changes) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 5245.3335  # H20
$ws.Cells.Item(20, 10).Value = 5600.3335  # J20
$ws.Cells.Item(20, 12).Value = 5600.3335  # L20
$ws.Cells.Item(20, 14).Value = -6094.3335  # N20
$ws.Cells.Item(35, 8).Value = 30000  # H35
$ws.Cells.Item(35, 10).Value = 30000  # J35
$ws.Cells.Item(35, 12).Value = 30000  # L35
$ws.Cells.Item(35, 14).Value = -30620  # N35
$ws.Cells.Item(82, 8).Value = 21967.23  # H82
$ws.Cells.Item(82, 10).Value = 24996  # J82
$ws.Cells.Item(82, 12).Value = 24996  # L82
$ws.Cells.Item(82, 14).Value = -25762  # N82
$ws.Cells.Item(85, 8).Value = 21967.23  # H85
$ws.Cells.Item(85, 10).Value = 24996  # J85
$ws.Cells.Item(85, 12).Value = 24996  # L85
$ws.Cells.Item(85, 14).Value = -27648  # N85
$ws.Cells.Item(86, 8).Value = 904.1081  # H86
$ws.Cells.Item(86, 9).Value = 897.2  # I86
$ws.Cells.Item(86, 10).Value = 1025  # J86
$ws.Cells.Item(86, 11).Value = 897.2  # K86
$ws.Cells.Item(86, 12).Value = 1025  # L86
$ws.Cells.Item(86, 13).Value = 225.8  # M86
$ws.Cells.Item(86, 14).Value = -3271  # N86
$ws.Cells.Item(89, 8).Value = 904.1081  # H89
$ws.Cells.Item(89, 9).Value = 897.2  # I89
$ws.Cells.Item(89, 10).Value = 1025  # J89
$ws.Cells.Item(89, 11).Value = 4486  # K89
$ws.Cells.Item(89, 12).Value = 5125  # L89
$ws.Cells.Item(89, 13).Value = 1130  # M89
$ws.Cells.Item(89, 14).Value = -16357  # N89
$ws.Cells.Item(94, 8).Value = 1231  # H94
$ws.Cells.Item(94, 9).Value = 1316.75  # I94
$ws.Cells.Item(94, 10).Value = 888  # J94
$ws.Cells.Item(94, 11).Value = 1316.75  # K94
$ws.Cells.Item(94, 12).Value = 888  # L94
$ws.Cells.Item(94, 13).Value = -865.75  # M94
$ws.Cells.Item(94, 14).Value = -1790  # N94
$ws.Cells.Item(97, 8).Value = 7471  # H97
$ws.Cells.Item(97, 9).Value = 7471  # I97
$ws.Cells.Item(97, 11).Value = 7471  # K97
$ws.Cells.Item(97, 13).Value = -6480  # M97
$ws.Cells.Item(99, 8).Value = 2174.0667  # H99
$ws.Cells.Item(99, 9).Value = 2511.182  # I99
$ws.Cells.Item(99, 11).Value = 2511.182  # K99
$ws.Cells.Item(99, 13).Value = -1013.182  # M99
$ws.Cells.Item(105, 8).Value = 4765686  # H105
$ws.Cells.Item(105, 9).Value = 5886772.5  # I105
$ws.Cells.Item(105, 10).Value = 1068.5  # J105
$ws.Cells.Item(105, 11).Value = 5886772.5  # K105
$ws.Cells.Item(105, 12).Value = 1068.5  # L105
$ws.Cells.Item(105, 13).Value = -5885025.5  # M105
$ws.Cells.Item(105, 14).Value = -4562.5  # N105
$ws.Cells.Item(107, 8).Value = 2681.111  # H107
$ws.Cells.Item(107, 9).Value = 2406.3845  # I107
$ws.Cells.Item(107, 10).Value = 3395.4  # J107
$ws.Cells.Item(107, 11).Value = 2406.3845  # K107
$ws.Cells.Item(107, 12).Value = 3395.4  # L107
$ws.Cells.Item(107, 13).Value = -486.3845000000001  # M107
$ws.Cells.Item(107, 14).Value = -7235.4  # N107

# ---- Sheet: CRP (76 cell changes) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4554.442  # H31
$ws.Cells.Item(31, 9).Value = 2938.25  # I31
$ws.Cells.Item(31, 10).Value = 7571.3335  # J31
$ws.Cells.Item(31, 11).Value = 2938.25  # K31
$ws.Cells.Item(31, 12).Value = 7571.3335  # L31
$ws.Cells.Item(31, 13).Value = -2643.25  # M31
$ws.Cells.Item(31, 14).Value = -8161.3335  # N31
$ws.Cells.Item(34, 8).Value = 4554.442  # H34
$ws.Cells.Item(34, 9).Value = 2938.25  # I34
$ws.Cells.Item(34, 10).Value = 7571.3335  # J34
$ws.Cells.Item(34, 11).Value = 2938.25  # K34
$ws.Cells.Item(34, 12).Value = 7571.3335  # L34
$ws.Cells.Item(34, 13).Value = -2736.25  # M34
$ws.Cells.Item(34, 14).Value = -7975.3335  # N34
$ws.Cells.Item(50, 8).Value = 14999  # H50
$ws.Cells.Item(50, 10).Value = 14999  # J50
$ws.Cells.Item(50, 12).Value = 14999  # L50
$ws.Cells.Item(50, 14).Value = -16249  # N50
$ws.Cells.Item(58, 8).Value = 1495.3334  # H58
$ws.Cells.Item(58, 9).Value = 1432.375  # I58
$ws.Cells.Item(58, 11).Value = 1432.375  # K58
$ws.Cells.Item(58, 13).Value = -1229.375  # M58
$ws.Cells.Item(59, 8).Value = 18494.8  # H59
$ws.Cells.Item(59, 10).Value = 18494.8  # J59
$ws.Cells.Item(59, 12).Value = 18494.8  # L59
$ws.Cells.Item(59, 14).Value = -20784.8  # N59
$ws.Cells.Item(68, 8).Value = 24998.75  # H68
$ws.Cells.Item(68, 10).Value = 24998.75  # J68
$ws.Cells.Item(68, 12).Value = 24998.75  # L68
$ws.Cells.Item(68, 14).Value = -26496.75  # N68
$ws.Cells.Item(71, 8).Value = 24998.75  # H71
$ws.Cells.Item(71, 10).Value = 24998.75  # J71
$ws.Cells.Item(71, 12).Value = 74996.25  # L71
$ws.Cells.Item(71, 14).Value = -82484.25  # N71
$ws.Cells.Item(74, 8).Value = 35950  # H74
$ws.Cells.Item(74, 10).Value = 36933.332  # J74
$ws.Cells.Item(74, 12).Value = 36933.332  # L74
$ws.Cells.Item(74, 14).Value = -38681.332  # N74
$ws.Cells.Item(77, 8).Value = 35950  # H77
$ws.Cells.Item(77, 10).Value = 36933.332  # J77
$ws.Cells.Item(77, 12).Value = 110799.996  # L77
$ws.Cells.Item(77, 14).Value = -119535.996  # N77
$ws.Cells.Item(86, 8).Value = 6291.6816  # H86
$ws.Cells.Item(86, 9).Value = 6174.231  # I86
$ws.Cells.Item(86, 10).Value = 6461.3335  # J86
$ws.Cells.Item(86, 11).Value = 6174.231  # K86
$ws.Cells.Item(86, 12).Value = 6461.3335  # L86
$ws.Cells.Item(86, 13).Value = -5051.231  # M86
$ws.Cells.Item(86, 14).Value = -8707.333500000001  # N86
$ws.Cells.Item(89, 8).Value = 6291.6816  # H89
$ws.Cells.Item(89, 9).Value = 6174.231  # I89
$ws.Cells.Item(89, 10).Value = 6461.3335  # J89
$ws.Cells.Item(89, 11).Value = 30871.155  # K89
$ws.Cells.Item(89, 12).Value = 32306.6675  # L89
$ws.Cells.Item(89, 13).Value = -25255.155  # M89
$ws.Cells.Item(89, 14).Value = -43538.6675  # N89
$ws.Cells.Item(94, 8).Value = 2598.375  # H94
$ws.Cells.Item(94, 9).Value = 3000.6  # I94
$ws.Cells.Item(94, 11).Value = 3000.6  # K94
$ws.Cells.Item(94, 13).Value = -2549.6  # M94
$ws.Cells.Item(99, 8).Value = 8817.272000000001  # H99
$ws.Cells.Item(99, 9).Value = 7998.8  # I99
$ws.Cells.Item(99, 11).Value = 7998.8  # K99
$ws.Cells.Item(99, 13).Value = -6500.8  # M99
$ws.Cells.Item(122, 8).Value = 1384.9375  # H122
$ws.Cells.Item(122, 9).Value = 1082.8572  # I122
$ws.Cells.Item(122, 11).Value = 3248.5716  # K122
$ws.Cells.Item(122, 13).Value = -798.5715999999998  # M122
$ws.Cells.Item(126, 8).Value = 8817.272000000001  # H126
$ws.Cells.Item(126, 9).Value = 7998.8  # I126
$ws.Cells.Item(126, 11).Value = 23996.4  # K126
$ws.Cells.Item(126, 13).Value = -21526.4  # M126
$ws.Cells.Item(136, 8).Value = 1495.3334  # H136
$ws.Cells.Item(136, 9).Value = 1432.375  # I136
$ws.Cells.Item(136, 11).Value = 4297.125  # K136
$ws.Cells.Item(136, 13).Value = -1747.125  # M136

# ---- Sheet: CUL (52 cell changes) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(10, 8).Value = 682.7368  # H10
$ws.Cells.Item(10, 9).Value = 256.77777  # I10
$ws.Cells.Item(10, 11).Value = 770.33331  # K10
$ws.Cells.Item(10, 13).Value = -631.33331  # M10
$ws.Cells.Item(17, 8).Value = 155  # H17
$ws.Cells.Item(17, 10).Value = 238  # J17
$ws.Cells.Item(17, 12).Value = 714  # L17
$ws.Cells.Item(17, 14).Value = -1052  # N17
$ws.Cells.Item(104, 8).Value = 170099.5  # H104
$ws.Cells.Item(104, 10).Value = 4119.6  # J104
$ws.Cells.Item(104, 12).Value = 12358.8  # L104
$ws.Cells.Item(104, 14).Value = -17600.8  # N104
$ws.Cells.Item(113, 8).Value = 625  # H113
$ws.Cells.Item(113, 10).Value = 500  # J113
$ws.Cells.Item(113, 12).Value = 1500  # L113
$ws.Cells.Item(113, 14).Value = -5840  # N113
$ws.Cells.Item(122, 8).Value = 735.6  # H122
$ws.Cells.Item(122, 9).Value = 192  # I122
$ws.Cells.Item(122, 10).Value = 819.2308  # J122
$ws.Cells.Item(122, 11).Value = 1728  # K122
$ws.Cells.Item(122, 12).Value = 7373.077200000001  # L122
$ws.Cells.Item(122, 13).Value = 722  # M122
$ws.Cells.Item(122, 14).Value = -12273.0772  # N122
$ws.Cells.Item(131, 8).Value = 104883.81  # H131
$ws.Cells.Item(131, 9).Value = 78146.69500000001  # I131
$ws.Cells.Item(131, 10).Value = 123177.63  # J131
$ws.Cells.Item(131, 11).Value = 234440.085  # K131
$ws.Cells.Item(131, 12).Value = 369532.89  # L131
$ws.Cells.Item(131, 13).Value = -229400.085  # M131
$ws.Cells.Item(131, 14).Value = -379612.89  # N131
$ws.Cells.Item(132, 8).Value = 2699.875  # H132
$ws.Cells.Item(132, 9).Value = 601  # I132
$ws.Cells.Item(132, 10).Value = 2999.7144  # J132
$ws.Cells.Item(132, 11).Value = 5409  # K132
$ws.Cells.Item(132, 12).Value = 26997.4296  # L132
$ws.Cells.Item(132, 13).Value = -2879  # M132
$ws.Cells.Item(132, 14).Value = -32057.4296  # N132
$ws.Cells.Item(133, 8).Value = 7490.3887  # H133
$ws.Cells.Item(133, 9).Value = 5916.2856  # I133
$ws.Cells.Item(133, 10).Value = 12999.75  # J133
$ws.Cells.Item(133, 11).Value = 17748.8568  # K133
$ws.Cells.Item(133, 12).Value = 38999.25  # L133
$ws.Cells.Item(133, 13).Value = -12688.8568  # M133
$ws.Cells.Item(133, 14).Value = -49119.25  # N133
$ws.Cells.Item(138, 8).Value = 3441.75  # H138
$ws.Cells.Item(138, 10).Value = 3666.3333  # J138
$ws.Cells.Item(138, 12).Value = 10998.9999  # L138
$ws.Cells.Item(138, 14).Value = -21278.9999  # N138
$ws.Cells.Item(140, 8).Value = 4379.05  # H140
$ws.Cells.Item(140, 9).Value = 3893.5  # I140
$ws.Cells.Item(140, 11).Value = 11680.5  # K140
$ws.Cells.Item(140, 13).Value = -6500.5  # M140

# ---- Sheet: GSM (19 cell changes) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(63, 8).Value = 21206  # H63
$ws.Cells.Item(63, 10).Value = 21557  # J63
$ws.Cells.Item(63, 12).Value = 21557  # L63
$ws.Cells.Item(63, 14).Value = -22929  # N63
$ws.Cells.Item(66, 8).Value = 21206  # H66
$ws.Cells.Item(66, 10).Value = 21557  # J66
$ws.Cells.Item(66, 12).Value = 64671  # L66
$ws.Cells.Item(66, 14).Value = -71535  # N66
$ws.Cells.Item(102, 8).Value = 3907.8  # H102
$ws.Cells.Item(102, 9).Value = 3366.05  # I102
$ws.Cells.Item(102, 11).Value = 3366.05  # K102
$ws.Cells.Item(102, 13).Value = -1744.05  # M102
$ws.Cells.Item(122, 8).Value = 2463.5186  # H122
$ws.Cells.Item(122, 9).Value = 2951.4119  # I122
$ws.Cells.Item(122, 10).Value = 1634.1  # J122
$ws.Cells.Item(122, 11).Value = 8854.235700000001  # K122
$ws.Cells.Item(122, 12).Value = 4902.299999999999  # L122
$ws.Cells.Item(122, 13).Value = -6404.235700000001  # M122
$ws.Cells.Item(122, 14).Value = -9802.299999999999  # N122

# ---- Sheet: LTW (57 cell changes) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 26463.715  # H7
$ws.Cells.Item(7, 9).Value = 38055  # I7
$ws.Cells.Item(7, 10).Value = 5599.4  # J7
$ws.Cells.Item(7, 11).Value = 38055  # K7
$ws.Cells.Item(7, 12).Value = 5599.4  # L7
$ws.Cells.Item(7, 13).Value = -37943  # M7
$ws.Cells.Item(7, 14).Value = -5823.4  # N7
$ws.Cells.Item(16, 8).Value = 4319.636  # H16
$ws.Cells.Item(16, 9).Value = 3279.8333  # I16
$ws.Cells.Item(16, 10).Value = 5567.4  # J16
$ws.Cells.Item(16, 11).Value = 3279.8333  # K16
$ws.Cells.Item(16, 12).Value = 5567.4  # L16
$ws.Cells.Item(16, 13).Value = -3109.8333  # M16
$ws.Cells.Item(16, 14).Value = -5907.4  # N16
$ws.Cells.Item(58, 8).Value = 10299.889  # H58
$ws.Cells.Item(58, 9).Value = 7999  # I58
$ws.Cells.Item(58, 11).Value = 7999  # K58
$ws.Cells.Item(58, 13).Value = -7739  # M58
$ws.Cells.Item(61, 8).Value = 13904644  # H61
$ws.Cells.Item(61, 9).Value = 16668705  # I61
$ws.Cells.Item(61, 10).Value = 84340  # J61
$ws.Cells.Item(61, 11).Value = 16668705  # K61
$ws.Cells.Item(61, 12).Value = 84340  # L61
$ws.Cells.Item(61, 13).Value = -16668503  # M61
$ws.Cells.Item(61, 14).Value = -84744  # N61
$ws.Cells.Item(100, 8).Value = 70081  # H100
$ws.Cells.Item(100, 9).Value = 4024.1667  # I100
$ws.Cells.Item(100, 10).Value = 169166.25  # J100
$ws.Cells.Item(100, 11).Value = 4024.1667  # K100
$ws.Cells.Item(100, 12).Value = 169166.25  # L100
$ws.Cells.Item(100, 13).Value = -3483.1667  # M100
$ws.Cells.Item(100, 14).Value = -170248.25  # N100
$ws.Cells.Item(113, 8).Value = 13904644  # H113
$ws.Cells.Item(113, 9).Value = 16668705  # I113
$ws.Cells.Item(113, 10).Value = 84340  # J113
$ws.Cells.Item(113, 11).Value = 16668705  # K113
$ws.Cells.Item(113, 12).Value = 84340  # L113
$ws.Cells.Item(113, 13).Value = -16666535  # M113
$ws.Cells.Item(113, 14).Value = -88680  # N113
$ws.Cells.Item(118, 8).Value = 0  # H118
$ws.Cells.Item(118, 10).Value = 0  # J118
$ws.Cells.Item(118, 12).Value = 0  # L118
$ws.Cells.Item(118, 14).ClearContents()  # N118: was -103314
$ws.Cells.Item(126, 8).Value = 26463.715  # H126
$ws.Cells.Item(126, 9).Value = 38055  # I126
$ws.Cells.Item(126, 10).Value = 5599.4  # J126
$ws.Cells.Item(126, 11).Value = 114165  # K126
$ws.Cells.Item(126, 12).Value = 16798.2  # L126
$ws.Cells.Item(126, 13).Value = -111695  # M126
$ws.Cells.Item(126, 14).Value = -21738.2  # N126
$ws.Cells.Item(132, 8).Value = 4358.45  # H132
$ws.Cells.Item(132, 9).Value = 4152  # I132
$ws.Cells.Item(132, 10).Value = 4977.8  # J132
$ws.Cells.Item(132, 11).Value = 12456  # K132
$ws.Cells.Item(132, 12).Value = 14933.4  # L132
$ws.Cells.Item(132, 13).Value = -9926  # M132
$ws.Cells.Item(132, 14).Value = -19993.4  # N132

# ---- Sheet: WVR (23 cell changes) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64, 8).Value = 25114  # H64
$ws.Cells.Item(64, 10).Value = 25114  # J64
$ws.Cells.Item(64, 12).Value = 25114  # L64
$ws.Cells.Item(64, 14).Value = -25610  # N64
$ws.Cells.Item(67, 8).Value = 25114  # H67
$ws.Cells.Item(67, 10).Value = 25114  # J67
$ws.Cells.Item(67, 12).Value = 25114  # L67
$ws.Cells.Item(67, 14).Value = -26830  # N67
$ws.Cells.Item(113, 8).Value = 4000  # H113
$ws.Cells.Item(113, 9).Value = 0  # I113
$ws.Cells.Item(113, 10).Value = 4000  # J113
$ws.Cells.Item(113, 11).Value = 0  # K113
$ws.Cells.Item(113, 12).Value = 12000  # L113
$ws.Cells.Item(113, 13).ClearContents()  # M113: was 225.5382999999999
$ws.Cells.Item(113, 14).Value = -16340  # N113
$ws.Cells.Item(122, 8).Value = 658.8946999999999  # H122
$ws.Cells.Item(122, 9).Value = 678.8889  # I122
$ws.Cells.Item(122, 11).Value = 2036.6667  # K122
$ws.Cells.Item(122, 13).Value = 413.3332999999998  # M122
$ws.Cells.Item(132, 8).Value = 14950.143  # H132
$ws.Cells.Item(132, 9).Value = 22487.75  # I132
$ws.Cells.Item(132, 11).Value = 67463.25  # K132
$ws.Cells.Item(132, 13).Value = -64933.25  # M132
